# Auto-generated edit script updating Sheets via scheduled runner
# Applies updated market-price-derived values (columns H-N) for specific
# Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1519.8518
$ws.Range("I98").Value = 1241.3636
$ws.Range("J98").Value = 2745.2
$ws.Range("K98").Value = 1241.3636
$ws.Range("L98").Value = 2745.2
$ws.Range("M98").Value = 256.6364000000001
$ws.Range("N98").Value = -5741.2

$ws.Range("H111").Value = 5308.619
$ws.Range("I111").Value = 3348.25
$ws.Range("J111").Value = 7922.4443
$ws.Range("K111").Value = 10044.75
$ws.Range("L111").Value = 23767.3329
$ws.Range("M111").Value = -6977.75
$ws.Range("N111").Value = -29901.3329

$ws.Range("H122").Value = 1519.8518
$ws.Range("I122").Value = 1241.3636
$ws.Range("J122").Value = 2745.2
$ws.Range("K122").Value = 3724.0908
$ws.Range("L122").Value = 8235.599999999999
$ws.Range("M122").Value = -1274.0908
$ws.Range("N122").Value = -13135.6

$ws.Range("H132").Value = 7832.3076
$ws.Range("I132").Value = 6985.923
$ws.Range("J132").Value = 9525.076999999999
$ws.Range("K132").Value = 20957.769
$ws.Range("L132").Value = 28575.231
$ws.Range("M132").Value = -18427.769
$ws.Range("N132").Value = -33635.231

$ws.Range("H138").Value = 2074.1843
$ws.Range("I138").Value = 2107.4736
$ws.Range("J138").Value = 2040.8948
$ws.Range("K138").Value = 6322.4208
$ws.Range("L138").Value = 6122.6844
$ws.Range("M138").Value = -1182.4208
$ws.Range("N138").Value = -16402.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4455614.5
$ws.Range("I32").Value = 6299.0605
$ws.Range("J32").Value = 37083930
$ws.Range("K32").Value = 6299.0605
$ws.Range("L32").Value = 37083930
$ws.Range("M32").Value = -6012.0605
$ws.Range("N32").Value = -37084504

$ws.Range("H37").Value = 18714.285
$ws.Range("J37").Value = 18714.285
$ws.Range("L37").Value = 18714.285
$ws.Range("N37").Value = -19260.285

$ws.Range("H132").Value = 2663520.5
$ws.Range("I132").Value = 1514.1724
$ws.Range("J132").Value = 6952308.5
$ws.Range("K132").Value = 4542.5172
$ws.Range("L132").Value = 20856925.5
$ws.Range("M132").Value = -2012.5172
$ws.Range("N132").Value = -20861985.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3110.104
$ws.Range("I134").Value = 1156.62
$ws.Range("J134").Value = 6727.6665
$ws.Range("K134").Value = 3469.86
$ws.Range("L134").Value = 20182.9995
$ws.Range("M134").Value = -934.8599999999997
$ws.Range("N134").Value = -25252.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = -9888

$ws.Range("H22").Value = 860.1818
$ws.Range("I22").Value = 276.4
$ws.Range("J22").Value = 1346.6666
$ws.Range("K22").Value = 276.4
$ws.Range("L22").Value = 1346.6666
$ws.Range("M22").Value = 73.60000000000002
$ws.Range("N22").Value = -2046.6666

$ws.Range("H74").Value = 39245.668
$ws.Range("J74").Value = 39245.668
$ws.Range("L74").Value = 39245.668
$ws.Range("N74").Value = -40993.668

$ws.Range("H77").Value = 39245.668
$ws.Range("J77").Value = 39245.668
$ws.Range("L77").Value = 117737.004
$ws.Range("N77").Value = -126473.004

$ws.Range("H94").Value = 1386.4
$ws.Range("J94").Value = 1398.8572
$ws.Range("L94").Value = 1398.8572
$ws.Range("N94").Value = -2300.8572

$ws.Range("H132").Value = 2886.9656
$ws.Range("I132").Value = 2314.4443
$ws.Range("J132").Value = 3144.6
$ws.Range("K132").Value = 6943.3329
$ws.Range("L132").Value = 9433.799999999999
$ws.Range("M132").Value = -4413.3329
$ws.Range("N132").Value = -14493.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 746396.9399999999
$ws.Range("I68").Value = 983
$ws.Range("J68").Value = 1058430.6
$ws.Range("K68").Value = 2949
$ws.Range("L68").Value = 3175291.8
$ws.Range("M68").Value = -2138
$ws.Range("N68").Value = -3176913.8

$ws.Range("H71").Value = 746396.9399999999
$ws.Range("I71").Value = 983
$ws.Range("J71").Value = 1058430.6
$ws.Range("K71").Value = 8847
$ws.Range("L71").Value = 9525875.4
$ws.Range("M71").Value = -4791
$ws.Range("N71").Value = -9533987.4

$ws.Range("H127").Value = 90909090
$ws.Range("J127").Value = 90909090
$ws.Range("L127").Value = 272727270
$ws.Range("N127").Value = -272737190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6370.8184
$ws.Range("I70").Value = 5072.375
$ws.Range("J70").Value = 9833.333000000001
$ws.Range("K70").Value = 5072.375
$ws.Range("L70").Value = 9833.333000000001
$ws.Range("M70").Value = -4802.375
$ws.Range("N70").Value = -10373.333

$ws.Range("H73").Value = 6370.8184
$ws.Range("I73").Value = 5072.375
$ws.Range("J73").Value = 9833.333000000001
$ws.Range("K73").Value = 5072.375
$ws.Range("L73").Value = 9833.333000000001
$ws.Range("M73").Value = -4136.375
$ws.Range("N73").Value = -11705.333

$ws.Range("H113").Value = 2040.125
$ws.Range("I113").Value = 2053.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2053.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 116.5
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 2422.3462
$ws.Range("I132").Value = 3168.125
$ws.Range("J132").Value = 2090.889
$ws.Range("K132").Value = 9504.375
$ws.Range("L132").Value = 6272.667
$ws.Range("M132").Value = -6974.375
$ws.Range("N132").Value = -11332.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3409.2
$ws.Range("I61").Value = 1716.2727
$ws.Range("J61").Value = 5478.3335
$ws.Range("K61").Value = 1716.2727
$ws.Range("L61").Value = 5478.3335
$ws.Range("M61").Value = -1514.2727
$ws.Range("N61").Value = -5882.3335

$ws.Range("H113").Value = 3409.2
$ws.Range("I113").Value = 1716.2727
$ws.Range("J113").Value = 5478.3335
$ws.Range("K113").Value = 1716.2727
$ws.Range("L113").Value = 5478.3335
$ws.Range("M113").Value = 453.7273
$ws.Range("N113").Value = -9818.333500000001

$ws.Range("H132").Value = 24392772
$ws.Range("I132").Value = 30305582
$ws.Range("J132").Value = 2424.125
$ws.Range("K132").Value = 90916746
$ws.Range("L132").Value = 7272.375
$ws.Range("M132").Value = -90914216
$ws.Range("N132").Value = -12332.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 781.35297
$ws.Range("I113").Value = 411.44446
$ws.Range("J113").Value = 914.52
$ws.Range("K113").Value = 1234.33338
$ws.Range("L113").Value = 2743.56
$ws.Range("M113").Value = 935.66662
$ws.Range("N113").Value = -7083.559999999999

$ws.Range("H132").Value = 2278.0513
$ws.Range("I132").Value = 1906.4
$ws.Range("J132").Value = 2669.2632
$ws.Range("K132").Value = 5719.200000000001
$ws.Range("L132").Value = 8007.7896
$ws.Range("M132").Value = -3189.200000000001
$ws.Range("N132").Value = -13067.7896
